$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8107.231
$ws.Range("I62").Value = 7491.75
$ws.Range("K62").Value = 7491.75
$ws.Range("M62").Value = -6867.75
$ws.Range("H65").Value = 8107.231
$ws.Range("I65").Value = 7491.75
$ws.Range("K65").Value = 37458.75
$ws.Range("M65").Value = -34338.75
$ws.Range("H70").Value = 18531376
$ws.Range("J70").Value = 24333.334
$ws.Range("L70").Value = 73000.00199999999
$ws.Range("N70").Value = -73540.00199999999
$ws.Range("H73").Value = 18531376
$ws.Range("J73").Value = 24333.334
$ws.Range("L73").Value = 73000.00199999999
$ws.Range("N73").Value = -74872.00199999999
$ws.Range("H81").Value = 400709
$ws.Range("J81").Value = 400709
$ws.Range("L81").Value = 400709
$ws.Range("N81").Value = -402705
$ws.Range("H84").Value = 400709
$ws.Range("J84").Value = 400709
$ws.Range("L84").Value = 1202127
$ws.Range("N84").Value = -1212111
$ws.Range("H106").Value = 2161.7334
$ws.Range("I106").Value = 2252.8572
$ws.Range("J106").Value = 886
$ws.Range("K106").Value = 2252.8572
$ws.Range("L106").Value = 886
$ws.Range("M106").Value = -1621.8572
$ws.Range("N106").Value = -2148
$ws.Range("H137").Value = 9857
$ws.Range("I137").Value = 1611.7858
$ws.Range("J137").Value = 14666.708
$ws.Range("K137").Value = 4835.357400000001
$ws.Range("L137").Value = 44000.124
$ws.Range("M137").Value = -2285.357400000001
$ws.Range("N137").Value = -49100.124
$ws.Range("H138").Value = 1419993.2
$ws.Range("I138").Value = 2478.375
$ws.Range("J138").Value = 2001537.9
$ws.Range("K138").Value = 7435.125
$ws.Range("L138").Value = 6004613.699999999
$ws.Range("M138").Value = -2295.125
$ws.Range("N138").Value = -6014893.699999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21946.314
$ws.Range("I32").Value = 25109.477
$ws.Range("K32").Value = 25109.477
$ws.Range("M32").Value = -24822.477
$ws.Range("H45").Value = 5783.6665
$ws.Range("I45").Value = 4942
$ws.Range("K45").Value = 4942
$ws.Range("M45").Value = -4565
$ws.Range("H61").Value = 8828.559999999999
$ws.Range("I61").Value = 5400.778
$ws.Range("K61").Value = 5400.778
$ws.Range("M61").Value = -5188.778
$ws.Range("H63").Value = 4131.773
$ws.Range("I63").Value = 2433.2222
$ws.Range("J63").Value = 5307.6924
$ws.Range("K63").Value = 2433.2222
$ws.Range("L63").Value = 5307.6924
$ws.Range("M63").Value = -1747.2222
$ws.Range("N63").Value = -6679.6924
$ws.Range("H66").Value = 4131.773
$ws.Range("I66").Value = 2433.2222
$ws.Range("J66").Value = 5307.6924
$ws.Range("K66").Value = 12166.111
$ws.Range("L66").Value = 26538.462
$ws.Range("M66").Value = -8734.111000000001
$ws.Range("N66").Value = -33402.462
$ws.Range("H132").Value = 1656.9231
$ws.Range("I132").Value = 1364.6666
$ws.Range("K132").Value = 4093.9998
$ws.Range("M132").Value = -1563.9998
$ws.Range("H136").Value = 8828.559999999999
$ws.Range("I136").Value = 5400.778
$ws.Range("K136").Value = 16202.334
$ws.Range("M136").Value = -13652.334
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 17520.273
$ws.Range("I82").Value = 10965.375
$ws.Range("K82").Value = 10965.375
$ws.Range("M82").Value = -10582.375
$ws.Range("H85").Value = 17520.273
$ws.Range("I85").Value = 10965.375
$ws.Range("K85").Value = 10965.375
$ws.Range("M85").Value = -9639.375
$ws.Range("H99").Value = 2363.5
$ws.Range("I99").Value = 1629.7142
$ws.Range("K99").Value = 1629.7142
$ws.Range("M99").Value = -131.7141999999999
$ws.Range("H134").Value = 11606
$ws.Range("I134").Value = 12096.275
$ws.Range("K134").Value = 36288.825
$ws.Range("M134").Value = -33753.825
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 7076998.5
$ws.Range("I6").Value = 10614498
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 10614498
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = -10614385
$ws.Range("N6").Value = -2226
$ws.Range("H31").Value = 2860380.2
$ws.Range("I31").Value = 20007400
$ws.Range("J31").Value = 2543.6
$ws.Range("K31").Value = 20007400
$ws.Range("L31").Value = 2543.6
$ws.Range("M31").Value = -20007105
$ws.Range("N31").Value = -3133.6
$ws.Range("H34").Value = 2860380.2
$ws.Range("I34").Value = 20007400
$ws.Range("J34").Value = 2543.6
$ws.Range("K34").Value = 20007400
$ws.Range("L34").Value = 2543.6
$ws.Range("M34").Value = -20007198
$ws.Range("N34").Value = -2947.6
$ws.Range("H132").Value = 2742.4443
$ws.Range("I132").Value = 2027.7693
$ws.Range("K132").Value = 6083.3079
$ws.Range("M132").Value = -3553.3079
$ws.Range("H134").Value = 2786.5356
$ws.Range("I134").Value = 2146.2273
$ws.Range("K134").Value = 6438.6819
$ws.Range("M134").Value = -3903.6819
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4800
$ws.Range("I56").Value = 4800
$ws.Range("K56").Value = 4800
$ws.Range("M56").Value = -4270
$ws.Range("H107").Value = 1734.5151
$ws.Range("I107").Value = 768.13336
$ws.Range("K107").Value = 2304.40008
$ws.Range("M107").Value = -384.4000800000003
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 30000
$ws.Range("J103").Value = 30000
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2014.3889
$ws.Range("I55").Value = 1219.2
$ws.Range("K55").Value = 1219.2
$ws.Range("M55").Value = -1046.2
$ws.Range("H93").Value = 1303
$ws.Range("I93").Value = 1293.75
$ws.Range("K93").Value = 1293.75
$ws.Range("M93").Value = -45.75
$ws.Range("H98").Value = 55000
$ws.Range("J98").Value = 55000
$ws.Range("L98").Value = 55000
$ws.Range("N98").Value = -60990
$ws.Range("H100").Value = 5000
$ws.Range("J100").Value = 5000
$ws.Range("L100").Value = 5000
$ws.Range("N100").Value = -6082
$ws.Range("H110").Value = 87855
$ws.Range("J110").Value = 87855
$ws.Range("L110").Value = 87855
$ws.Range("N110").Value = -96035
$ws.Range("H132").Value = 6748.6
$ws.Range("J132").Value = 5435.875
$ws.Range("L132").Value = 16307.625
$ws.Range("N132").Value = -21367.625
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10000
$ws.Range("I62").Value = 10000
$ws.Range("K62").Value = 10000
$ws.Range("M62").Value = -9376
$ws.Range("H65").Value = 10000
$ws.Range("I65").Value = 10000
$ws.Range("K65").Value = 50000
$ws.Range("M65").Value = -46880
$ws.Range("H113").Value = 745.8919
$ws.Range("I113").Value = 789.76666
$ws.Range("J113").Value = 557.8570999999999
$ws.Range("K113").Value = 2369.29998
$ws.Range("L113").Value = 1673.5713
$ws.Range("M113").Value = -199.2999799999998
$ws.Range("N113").Value = -6013.5713
$ws.Range("H122").Value = 3238.1724
$ws.Range("I122").Value = 3452.4285
$ws.Range("J122").Value = 2675.75
$ws.Range("K122").Value = 10357.2855
$ws.Range("L122").Value = 8027.25
$ws.Range("M122").Value = -7907.2855
$ws.Range("N122").Value = -12927.25
$ws.Range("H136").Value = 7349.531
$ws.Range("I136").Value = 8246.362999999999
$ws.Range("K136").Value = 24739.089
$ws.Range("M136").Value = -22189.089
